$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -match $pattern) {
            return $idx
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Merge the two runs that straddle the "_GoBack" bookmark in the intro
#    paragraph into a single run, and drop the bookmark. A Find/Replace
#    across the bookmark boundary (replacing text with itself) makes Word
#    recombine the run and removes the now-redundant bookmark.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$introOld = "builds on the Open"
$find.Execute($introOld, $true, $false, $false, $false, $false, $true, 1, $false, $introOld, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Turn the blank paragraph under "The following items will be included
#    in the scope of work:" into a bulleted list item (Inclusions list).
# ---------------------------------------------------------------------------
$exclIntroIdx = Find-ParagraphIndex("following items will not be included")
$pExclIntro = $d.Paragraphs($exclIntroIdx)
$pExclIntro.Range.InsertParagraphAfter()

$pExclBullet = $d.Paragraphs($exclIntroIdx + 1)
$pExclBullet.Range.Text = "Messages sent from server to client"
$pExclBullet.Style = "List Paragraph"
$pExclBullet.Range.ListFormat.ApplyBulletDefault()

$inclIntroIdx = Find-ParagraphIndex("following items will be included")
$pInclBlank = $d.Paragraphs($inclIntroIdx + 1)
$pInclBlank.Range.Text = "Messages sent from client to server"
$pInclBlank.Style = "List Paragraph"
$pInclBlank.Range.ListFormat.ApplyBulletDefault()

# ---------------------------------------------------------------------------
# 3) Schedule ("Hours Breakdown") table edits:
#      - "Test Pit" row gains a lastRenderedPageBreak marker on its run.
#      - "Acceptance Testing" row's "2 days" becomes two runs "1" + " days".
#      - "TOTAL HOURS" row loses its lastRenderedPageBreak marker.
#      - "TOTAL HOURS" total changes from "18" to "17" (kept as two runs).
# ---------------------------------------------------------------------------
$hoursTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables($i)
    if ($tbl.Cell(1,1).Range.Text -match "Work Item") {
        $hoursTable = $tbl
    }
}

$testPitRow = 0
$acceptRow = 0
$totalRow = 0
for ($r = 1; $r -le $hoursTable.Rows.Count; $r++) {
    $label = $hoursTable.Cell($r,1).Range.Text
    if ($label -match "^Test Pit") { $testPitRow = $r }
    if ($label -match "^Acceptance Testing") { $acceptRow = $r }
    if ($label -match "^TOTAL HOURS") { $totalRow = $r }
}

$pTestPit = $hoursTable.Cell($testPitRow,1).Range.Paragraphs(1)
$pTestPit.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Test Pit</w:t></w:r></w:p>')

$pAccept = $hoursTable.Cell($acceptRow,2).Range.Paragraphs(1)
$pAccept.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space="preserve"> days</w:t></w:r></w:p>')

$pTotalLabel = $hoursTable.Cell($totalRow,1).Range.Paragraphs(1)
$pTotalLabel.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>TOTAL HOURS</w:t></w:r></w:p>')

$pTotalValue = $hoursTable.Cell($totalRow,2).Range.Paragraphs(1)
$pTotalValue.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>17</w:t></w:r><w:r><w:t xml:space="preserve"> days</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 4) "Test Environment (Huawei)" deliverable paragraph: add a hanging
#    indent and move the "_GoBack" bookmark here (it was removed from the
#    intro paragraph in step 1, so it is free to be reassigned id 0).
# ---------------------------------------------------------------------------
$testEnvIdx = Find-ParagraphIndex("Huawei will provide a work test environment")
$pTestEnv = $d.Paragraphs($testEnvIdx)
$testEnvXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r><w:t>Huawei will provide a work test environment for validation of the fuzzing definition. If the protocol is supported, Huawei will provide a configuration for the Deja vu Security''s lab containing two Huawei AR series routers.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pTestEnv.Range.InsertXML($testEnvXml)

Write-Output "done"
